$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so values like "19.60" or "0.190"
# are not auto-converted to numbers, losing trailing zeros / thousand separators.
$textCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "E11", "E12", "D13", "E13", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($c in $textCells) {
  $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.886.58'
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").Value = '2.988.88'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '542.83'
$ws.Range("E5").Value = '  -3.27%  '
$ws.Range("D6").Value = '151.46'
$ws.Range("E6").Value = '  -4.05%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").Value = '3.005.51'
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("E11").Value = '  -5.65%  '
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").Value = '3.518.88'
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '61.940.00'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '23.91'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '3.002.33'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").Value = '0.0000147'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").Value = '5.17'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '12.03'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").Value = '378.67'
$ws.Range("E21").Value = '  -4.77%  '
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '5.67'
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '66.06'
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.118.53'
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = '0.469'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '0.190'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0927'
$ws.Range("E30").Value = '  -5.86%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '8.27'
$ws.Range("E31").Value = '  -6.24%  '
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '1.73'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '20.49'
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '160.95'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.60'
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '5.91'
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.27'
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '1.55'
$ws.Range("E40").Value = '  -4.15%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '37.54'
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.416.23'
$ws.Range("E42").Value = '  -4.65%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.90'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '22.04'
$ws.Range("E44").Value = '  -4.04%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.673'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '0.0591'
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '5.36'
$ws.Range("E47").Value = '  +6.75%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").Value = '0.997'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0244'
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = '268.90'
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.0951'
$ws.Range("E51").Value = '  -0.31%  '
